$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documents")

$ws.Range("I3").Hyperlinks.Delete()

$ws.Range("I3").Value = "Z:\USERS\Concept design and Options Report_rev2_25.7 Including Appendices.pdf"
$ws.Range("I3").Style = "Hyperlink"
$ws.Range("I6").Value = "Z:\USERS\Concept design and Options Report_rev2_25.7 Including Appendices.pdf"
$ws.Range("I6").Style = "Hyperlink"

$u = "https://dubaiholding-my.sharepoint.com/:b:/r/personal/arun_naidu_dhre_ae/Documents/Documents/NME/MARINE/DMS%20148763%20-%20Marine%20Works%20Including%20Dredging%20and%20Land%20Reclamation%20(Jan%20De%20Nul)/EMPLOYERS%20REPRESENTATIVES%20INSTRUCTION/ERI%2301/NKL-LT-12161%20-%20ERI%2301%20-%20Additional%20Sand%20Stockpiles.pdf?csf=1&web=1&e=7f8NQc"
$ws.Hyperlinks.Add($ws.Range("I2"), $u)
$ws.Hyperlinks.Add($ws.Range("I4"), $u)
$ws.Hyperlinks.Add($ws.Range("I5"), $u)
$ws.Hyperlinks.Add($ws.Range("I7"), $u)

# new row 8
$ws.Range("A8").Value = "PC1.0"
$ws.Range("B8").Value = "DMS 149600"
$ws.Range("C8").Value = "PC#01 - Payment Certificate.pdf"
$ws.Range("D8").Value = "PAYMENT CERTIFICATE"
$ws.Range("E8").Value = 45827
$ws.Range("E8").Style = $ws.Range("E7").Style
$ws.Range("F8").Value = "PC # 01.0"
$ws.Range("G8").Value = "Payment Certificate # 01"
$ws.Range("H8").Value = "PC-01.PDF"
$ws.Range("I8").Value = "Payment Certificate # 01"
$ws.Range("J8").Value = 45827
$ws.Range("J8").Style = $ws.Range("J7").Style
$ws.Range("K8").Value = 0

$u2 = "https://dubaiholding-my.sharepoint.com/:b:/g/personal/arun_naidu_dhre_ae/EUysg7T3ajdIsj6J81h2b_8BqMv40tjbbnmHigxS5cyl8w?email=Arun.Naidu%40dhre.ae&e=fNZcxb"
$ws.Hyperlinks.Add($ws.Range("I8"), $u2)
$ws.Range("I8").Style = "Hyperlink"
